$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 18.1667185
$ws.Cells.Item(2, 8).Value = 36.333437
$ws.Cells.Item(2, 9).Value = 0.1223713430663405
$ws.Cells.Item(2, 10).Value = 0.09065557961739065
$ws.Cells.Item(2, 13).Value = 50.6072485
$ws.Cells.Item(2, 14).Value = 101.214497
$ws.Cells.Item(2, 15).Value = 0.2036920096625967
$ws.Cells.Item(2, 16).Value = 0.1535170070198019
$ws.Cells.Item(2, 17).Value = 919.3676375590472
$ws.Cells.Item(2, 18).Value = 3677.470550236189
$ws.Cells.Item(2, 19).Value = 0.02492606479429397
$ws.Cells.Item(2, 20).Value = 0.01391717325250717
$ws.Cells.Item(3, 7).Value = 18.1667185
$ws.Cells.Item(3, 8).Value = 36.333437
$ws.Cells.Item(3, 9).Value = 0.1223713430663405
$ws.Cells.Item(3, 10).Value = 0.09065557961739065
$ws.Cells.Item(3, 15).Value = 0.5188519664463093
$ws.Cells.Item(3, 16).Value = 0.5865664620849566
$ws.Cells.Item(3, 17).Value = 2341.847907656057
$ws.Cells.Item(3, 18).Value = 14051.08744593634
$ws.Cells.Item(3, 19).Value = 0.0634926119866467
$ws.Cells.Item(3, 20).Value = 0.05317552260443394
$ws.Cells.Item(4, 7).Value = 18.1667185
$ws.Cells.Item(4, 8).Value = 36.333437
$ws.Cells.Item(4, 9).Value = 0.1223713430663405
$ws.Cells.Item(4, 10).Value = 0.09065557961739065
$ws.Cells.Item(4, 13).Value = 16.762851
$ws.Cells.Item(4, 14).Value = 50.288553
$ws.Cells.Item(4, 15).Value = 0.06746975797083039
$ws.Cells.Item(4, 16).Value = 0.07627512236628199
$ws.Cells.Item(4, 17).Value = 304.5259953744435
$ws.Cells.Item(4, 18).Value = 1827.155972246661
$ws.Cells.Item(4, 19).Value = 0.008256364899251448
$ws.Cells.Item(4, 20).Value = 0.006914765428502691
$ws.Cells.Item(5, 7).Value = 18.1667185
$ws.Cells.Item(5, 8).Value = 36.333437
$ws.Cells.Item(5, 9).Value = 0.1223713430663405
$ws.Cells.Item(5, 10).Value = 0.09065557961739065
$ws.Cells.Item(5, 13).Value = 35.4375075
$ws.Cells.Item(5, 14).Value = 70.87501499999999
$ws.Cells.Item(5, 15).Value = 0.1426344512705199
$ws.Cells.Item(5, 16).Value = 0.1074996220678108
$ws.Cells.Item(5, 17).Value = 643.7832230941387
$ws.Cells.Item(5, 18).Value = 2575.132892376555
$ws.Cells.Item(5, 19).Value = 0.01745436936950402
$ws.Cells.Item(5, 20).Value = 0.009745440547207824
$ws.Cells.Item(6, 7).Value = 18.1667185
$ws.Cells.Item(6, 8).Value = 36.333437
$ws.Cells.Item(6, 9).Value = 0.1223713430663405
$ws.Cells.Item(6, 10).Value = 0.09065557961739065
$ws.Cells.Item(6, 13).Value = 3.403012
$ws.Cells.Item(6, 14).Value = 10.209036
$ws.Cells.Item(6, 15).Value = 0.01369697768069593
$ws.Cells.Item(6, 16).Value = 0.01548454715214769
$ws.Cells.Item(6, 17).Value = 61.82156105612201
$ws.Cells.Item(6, 18).Value = 370.9293663367321
$ws.Cells.Item(6, 19).Value = 0.001676117554736451
$ws.Cells.Item(6, 20).Value = 0.001403760597190764
$ws.Cells.Item(7, 7).Value = 18.1667185
$ws.Cells.Item(7, 8).Value = 36.333437
$ws.Cells.Item(7, 9).Value = 0.1223713430663405
$ws.Cells.Item(7, 10).Value = 0.09065557961739065
$ws.Cells.Item(7, 13).Value = 13.330536
$ws.Cells.Item(7, 14).Value = 39.991608
$ws.Cells.Item(7, 15).Value = 0.05365483696904789
$ws.Cells.Item(7, 16).Value = 0.06065723930900103
$ws.Cells.Item(7, 17).Value = 242.172094966116
$ws.Cells.Item(7, 18).Value = 1453.032569796696
$ws.Cells.Item(7, 19).Value = 0.006565814461907929
$ws.Cells.Item(7, 20).Value = 0.005498917187548262
$ws.Cells.Item(8, 9).Value = 0.6434767683046462
$ws.Cells.Item(8, 10).Value = 0.715054169619088
$ws.Cells.Item(8, 13).Value = 50.6072485
$ws.Cells.Item(8, 14).Value = 101.214497
$ws.Cells.Item(8, 15).Value = 0.2036920096625967
$ws.Cells.Item(8, 16).Value = 0.1535170070198019
$ws.Cells.Item(8, 17).Value = 4834.39751069543
$ws.Cells.Item(8, 18).Value = 29006.38506417258
$ws.Cells.Item(8, 19).Value = 0.1310710761071665
$ws.Cells.Item(8, 20).Value = 0.1097729759769522
$ws.Cells.Item(9, 9).Value = 0.6434767683046462
$ws.Cells.Item(9, 10).Value = 0.715054169619088
$ws.Cells.Item(9, 15).Value = 0.5188519664463093
$ws.Cells.Item(9, 16).Value = 0.5865664620849566
$ws.Cells.Item(9, 19).Value = 0.3338691865973818
$ws.Cells.Item(9, 20).Value = 0.4194267944725649
$ws.Cells.Item(10, 9).Value = 0.6434767683046462
$ws.Cells.Item(10, 10).Value = 0.715054169619088
$ws.Cells.Item(10, 13).Value = 16.762851
$ws.Cells.Item(10, 14).Value = 50.288553
$ws.Cells.Item(10, 15).Value = 0.06746975797083039
$ws.Cells.Item(10, 16).Value = 0.07627512236628199
$ws.Cells.Item(10, 17).Value = 1601.317746934185
$ws.Cells.Item(10, 18).Value = 14411.85972240767
$ws.Cells.Item(10, 19).Value = 0.04341522181736658
$ws.Cells.Item(10, 20).Value = 0.05454084428621609
$ws.Cells.Item(11, 9).Value = 0.6434767683046462
$ws.Cells.Item(11, 10).Value = 0.715054169619088
$ws.Cells.Item(11, 13).Value = 35.4375075
$ws.Cells.Item(11, 14).Value = 70.87501499999999
$ws.Cells.Item(11, 15).Value = 0.1426344512705199
$ws.Cells.Item(11, 16).Value = 0.1074996220678108
$ws.Cells.Item(11, 17).Value = 3385.266006770762
$ws.Cells.Item(11, 18).Value = 20311.59604062457
$ws.Cells.Item(11, 19).Value = 0.0917819557524607
$ws.Cells.Item(11, 20).Value = 0.0768680529920642
$ws.Cells.Item(12, 9).Value = 0.6434767683046462
$ws.Cells.Item(12, 10).Value = 0.715054169619088
$ws.Cells.Item(12, 13).Value = 3.403012
$ws.Cells.Item(12, 14).Value = 10.209036
$ws.Cells.Item(12, 15).Value = 0.01369697768069593
$ws.Cells.Item(12, 16).Value = 0.01548454715214769
$ws.Cells.Item(12, 17).Value = 325.0821419715534
$ws.Cells.Item(12, 18).Value = 2925.73927774398
$ws.Cells.Item(12, 19).Value = 0.008813686933515086
$ws.Cells.Item(12, 20).Value = 0.01107229000580658
$ws.Cells.Item(13, 9).Value = 0.6434767683046462
$ws.Cells.Item(13, 10).Value = 0.715054169619088
$ws.Cells.Item(13, 13).Value = 13.330536
$ws.Cells.Item(13, 14).Value = 39.991608
$ws.Cells.Item(13, 15).Value = 0.05365483696904789
$ws.Cells.Item(13, 16).Value = 0.06065723930900103
$ws.Cells.Item(13, 17).Value = 1273.43635476716
$ws.Cells.Item(13, 18).Value = 11460.92719290444
$ws.Cells.Item(13, 19).Value = 0.0345256410967556
$ws.Cells.Item(13, 20).Value = 0.04337321188548404
$ws.Cells.Item(14, 7).Value = 7.527206666666667
$ws.Cells.Item(14, 8).Value = 22.58162
$ws.Cells.Item(14, 9).Value = 0.05070339969972629
$ws.Cells.Item(14, 10).Value = 0.05634341308805058
$ws.Cells.Item(14, 13).Value = 50.6072485
$ws.Cells.Item(14, 14).Value = 101.214497
$ws.Cells.Item(14, 15).Value = 0.2036920096625967
$ws.Cells.Item(14, 16).Value = 0.1535170070198019
$ws.Cells.Item(14, 17).Value = 380.9312182908567
$ws.Cells.Item(14, 18).Value = 2285.58730974514
$ws.Cells.Item(14, 19).Value = 0.01032787738156315
$ws.Cells.Item(14, 20).Value = 0.00864967214255786
$ws.Cells.Item(15, 7).Value = 7.527206666666667
$ws.Cells.Item(15, 8).Value = 22.58162
$ws.Cells.Item(15, 9).Value = 0.05070339969972629
$ws.Cells.Item(15, 10).Value = 0.05634341308805058
$ws.Cells.Item(15, 15).Value = 0.5188519664463093
$ws.Cells.Item(15, 16).Value = 0.5865664620849566
$ws.Cells.Item(15, 17).Value = 970.3223607955423
$ws.Cells.Item(15, 18).Value = 8732.901247159882
$ws.Cells.Item(15, 19).Value = 0.02630755863971619
$ws.Cells.Item(15, 20).Value = 0.03304915647684906
$ws.Cells.Item(16, 7).Value = 7.527206666666667
$ws.Cells.Item(16, 8).Value = 22.58162
$ws.Cells.Item(16, 9).Value = 0.05070339969972629
$ws.Cells.Item(16, 10).Value = 0.05634341308805058
$ws.Cells.Item(16, 13).Value = 16.762851
$ws.Cells.Item(16, 14).Value = 50.288553
$ws.Cells.Item(16, 15).Value = 0.06746975797083039
$ws.Cells.Item(16, 16).Value = 0.07627512236628199
$ws.Cells.Item(16, 17).Value = 126.17744379954
$ws.Cells.Item(16, 18).Value = 1135.59699419586
$ws.Cells.Item(16, 19).Value = 0.003420946106038807
$ws.Cells.Item(16, 20).Value = 0.004297600727825032
$ws.Cells.Item(17, 7).Value = 7.527206666666667
$ws.Cells.Item(17, 8).Value = 22.58162
$ws.Cells.Item(17, 9).Value = 0.05070339969972629
$ws.Cells.Item(17, 10).Value = 0.05634341308805058
$ws.Cells.Item(17, 13).Value = 35.4375075
$ws.Cells.Item(17, 14).Value = 70.87501499999999
$ws.Cells.Item(17, 15).Value = 0.1426344512705199
$ws.Cells.Item(17, 16).Value = 0.1074996220678108
$ws.Cells.Item(17, 17).Value = 266.74544270405
$ws.Cells.Item(17, 18).Value = 1600.4726562243
$ws.Cells.Item(17, 19).Value = 0.007232051593720304
$ws.Cells.Item(17, 20).Value = 0.00605689561297598
$ws.Cells.Item(18, 7).Value = 7.527206666666667
$ws.Cells.Item(18, 8).Value = 22.58162
$ws.Cells.Item(18, 9).Value = 0.05070339969972629
$ws.Cells.Item(18, 10).Value = 0.05634341308805058
$ws.Cells.Item(18, 13).Value = 3.403012
$ws.Cells.Item(18, 14).Value = 10.209036
$ws.Cells.Item(18, 15).Value = 0.01369697768069593
$ws.Cells.Item(18, 16).Value = 0.01548454715214769
$ws.Cells.Item(18, 17).Value = 25.61517461314667
$ws.Cells.Item(18, 18).Value = 230.53657151832
$ws.Cells.Item(18, 19).Value = 0.0006944833340225557
$ws.Cells.Item(18, 20).Value = 0.0008724522366748544
$ws.Cells.Item(19, 7).Value = 7.527206666666667
$ws.Cells.Item(19, 8).Value = 22.58162
$ws.Cells.Item(19, 9).Value = 0.05070339969972629
$ws.Cells.Item(19, 10).Value = 0.05634341308805058
$ws.Cells.Item(19, 13).Value = 13.330536
$ws.Cells.Item(19, 14).Value = 39.991608
$ws.Cells.Item(19, 15).Value = 0.05365483696904789
$ws.Cells.Item(19, 16).Value = 0.06065723930900103
$ws.Cells.Item(19, 17).Value = 100.34169944944
$ws.Cells.Item(19, 18).Value = 903.0752950449601
$ws.Cells.Item(19, 19).Value = 0.002720482644665286
$ws.Cells.Item(19, 20).Value = 0.003417635891167785
$ws.Cells.Item(20, 7).Value = 26.4148145
$ws.Cells.Item(20, 8).Value = 52.829629
$ws.Cells.Item(20, 9).Value = 0.1779306663013051
$ws.Cells.Item(20, 10).Value = 0.1318152377923044
$ws.Cells.Item(20, 13).Value = 50.6072485
$ws.Cells.Item(20, 14).Value = 101.214497
$ws.Cells.Item(20, 15).Value = 0.2036920096625967
$ws.Cells.Item(20, 16).Value = 0.1535170070198019
$ws.Cells.Item(20, 17).Value = 1336.781081482903
$ws.Cells.Item(20, 18).Value = 5347.124325931612
$ws.Cells.Item(20, 19).Value = 0.0362430549995177
$ws.Cells.Item(20, 20).Value = 0.02023588078547805
$ws.Cells.Item(21, 7).Value = 26.4148145
$ws.Cells.Item(21, 8).Value = 52.829629
$ws.Cells.Item(21, 9).Value = 0.1779306663013051
$ws.Cells.Item(21, 10).Value = 0.1318152377923044
$ws.Cells.Item(21, 15).Value = 0.5188519664463093
$ws.Cells.Item(21, 16).Value = 0.5865664620849566
$ws.Cells.Item(21, 17).Value = 3405.099169007758
$ws.Cells.Item(21, 18).Value = 20430.59501404655
$ws.Cells.Item(21, 19).Value = 0.09231967610153417
$ws.Cells.Item(21, 20).Value = 0.07731839768071923
$ws.Cells.Item(22, 7).Value = 26.4148145
$ws.Cells.Item(22, 8).Value = 52.829629
$ws.Cells.Item(22, 9).Value = 0.1779306663013051
$ws.Cells.Item(22, 10).Value = 0.1318152377923044
$ws.Cells.Item(22, 13).Value = 16.762851
$ws.Cells.Item(22, 14).Value = 50.288553
$ws.Cells.Item(22, 15).Value = 0.06746975797083039
$ws.Cells.Item(22, 16).Value = 0.07627512236628199
$ws.Cells.Item(22, 17).Value = 442.7875996561395
$ws.Cells.Item(22, 18).Value = 2656.725597936837
$ws.Cells.Item(22, 19).Value = 0.01200493899093764
$ws.Cells.Item(22, 20).Value = 0.01005422339234857
$ws.Cells.Item(23, 7).Value = 26.4148145
$ws.Cells.Item(23, 8).Value = 52.829629
$ws.Cells.Item(23, 9).Value = 0.1779306663013051
$ws.Cells.Item(23, 10).Value = 0.1318152377923044
$ws.Cells.Item(23, 13).Value = 35.4375075
$ws.Cells.Item(23, 14).Value = 70.87501499999999
$ws.Cells.Item(23, 15).Value = 0.1426344512705199
$ws.Cells.Item(23, 16).Value = 0.1074996220678108
$ws.Cells.Item(23, 17).Value = 936.0751869548586
$ws.Cells.Item(23, 18).Value = 3744.300747819434
$ws.Cells.Item(23, 19).Value = 0.02537904295208464
$ws.Cells.Item(23, 20).Value = 0.01417008824545133
$ws.Cells.Item(24, 7).Value = 26.4148145
$ws.Cells.Item(24, 8).Value = 52.829629
$ws.Cells.Item(24, 9).Value = 0.1779306663013051
$ws.Cells.Item(24, 10).Value = 0.1318152377923044
$ws.Cells.Item(24, 13).Value = 3.403012
$ws.Cells.Item(24, 14).Value = 10.209036
$ws.Cells.Item(24, 15).Value = 0.01369697768069593
$ws.Cells.Item(24, 16).Value = 0.01548454715214769
$ws.Cells.Item(24, 17).Value = 89.889930721274
$ws.Cells.Item(24, 18).Value = 539.339584327644
$ws.Cells.Item(24, 19).Value = 0.002437112365040331
$ws.Cells.Item(24, 20).Value = 0.002041099264966497
$ws.Cells.Item(25, 7).Value = 26.4148145
$ws.Cells.Item(25, 8).Value = 52.829629
$ws.Cells.Item(25, 9).Value = 0.1779306663013051
$ws.Cells.Item(25, 10).Value = 0.1318152377923044
$ws.Cells.Item(25, 13).Value = 13.330536
$ws.Cells.Item(25, 14).Value = 39.991608
$ws.Cells.Item(25, 15).Value = 0.05365483696904789
$ws.Cells.Item(25, 16).Value = 0.06065723930900103
$ws.Cells.Item(25, 17).Value = 352.123635625572
$ws.Cells.Item(25, 18).Value = 2112.741813753432
$ws.Cells.Item(25, 19).Value = 0.009546840892190586
$ws.Cells.Item(25, 20).Value = 0.007995548423340684
$ws.Cells.Item(26, 7).Value = 0.5285683333333333
$ws.Cells.Item(26, 8).Value = 1.585705
$ws.Cells.Item(26, 9).Value = 0.003560445814819949
$ws.Cells.Item(26, 10).Value = 0.00395649346020291
$ws.Cells.Item(26, 13).Value = 50.6072485
$ws.Cells.Item(26, 14).Value = 101.214497
$ws.Cells.Item(26, 15).Value = 0.2036920096625967
$ws.Cells.Item(26, 16).Value = 0.1535170070198019
$ws.Cells.Item(26, 17).Value = 26.74938899423083
$ws.Cells.Item(26, 18).Value = 160.496333965385
$ws.Cells.Item(26, 19).Value = 0.0007252343633154571
$ws.Cells.Item(26, 20).Value = 0.0006073890343037705
$ws.Cells.Item(27, 7).Value = 0.5285683333333333
$ws.Cells.Item(27, 8).Value = 1.585705
$ws.Cells.Item(27, 9).Value = 0.003560445814819949
$ws.Cells.Item(27, 10).Value = 0.00395649346020291
$ws.Cells.Item(27, 15).Value = 0.5188519664463093
$ws.Cells.Item(27, 16).Value = 0.5865664620849566
$ws.Cells.Item(27, 17).Value = 68.13705213024112
$ws.Cells.Item(27, 18).Value = 613.23346917217
$ws.Cells.Item(27, 19).Value = 0.001847344312444862
$ws.Cells.Item(27, 20).Value = 0.002320746371213489
$ws.Cells.Item(28, 7).Value = 0.5285683333333333
$ws.Cells.Item(28, 8).Value = 1.585705
$ws.Cells.Item(28, 9).Value = 0.003560445814819949
$ws.Cells.Item(28, 10).Value = 0.00395649346020291
$ws.Cells.Item(28, 13).Value = 16.762851
$ws.Cells.Item(28, 14).Value = 50.288553
$ws.Cells.Item(28, 15).Value = 0.06746975797083039
$ws.Cells.Item(28, 16).Value = 0.07627512236628199
$ws.Cells.Item(28, 17).Value = 8.860312214985001
$ws.Cells.Item(28, 18).Value = 79.74280993486499
$ws.Cells.Item(28, 19).Value = 0.000240222417394158
$ws.Cells.Item(28, 20).Value = 0.0003017820228183714
$ws.Cells.Item(29, 7).Value = 0.5285683333333333
$ws.Cells.Item(29, 8).Value = 1.585705
$ws.Cells.Item(29, 9).Value = 0.003560445814819949
$ws.Cells.Item(29, 10).Value = 0.00395649346020291
$ws.Cells.Item(29, 13).Value = 35.4375075
$ws.Cells.Item(29, 14).Value = 70.87501499999999
$ws.Cells.Item(29, 15).Value = 0.1426344512705199
$ws.Cells.Item(29, 16).Value = 0.1074996220678108
$ws.Cells.Item(29, 17).Value = 18.7311442767625
$ws.Cells.Item(29, 18).Value = 112.386865660575
$ws.Cells.Item(29, 19).Value = 0.0005078422350752626
$ws.Cells.Item(29, 20).Value = 0.0004253215516855777
$ws.Cells.Item(30, 7).Value = 0.5285683333333333
$ws.Cells.Item(30, 8).Value = 1.585705
$ws.Cells.Item(30, 9).Value = 0.003560445814819949
$ws.Cells.Item(30, 10).Value = 0.00395649346020291
$ws.Cells.Item(30, 13).Value = 3.403012
$ws.Cells.Item(30, 14).Value = 10.209036
$ws.Cells.Item(30, 15).Value = 0.01369697768069593
$ws.Cells.Item(30, 16).Value = 0.01548454715214769
$ws.Cells.Item(30, 17).Value = 1.798724381153333
$ws.Cells.Item(30, 18).Value = 16.18851943038
$ws.Cells.Item(30, 19).Value = 0.00004876734685891608
$ws.Cells.Item(30, 20).Value = 0.00006126450954167593
$ws.Cells.Item(31, 7).Value = 0.5285683333333333
$ws.Cells.Item(31, 8).Value = 1.585705
$ws.Cells.Item(31, 9).Value = 0.003560445814819949
$ws.Cells.Item(31, 10).Value = 0.00395649346020291
$ws.Cells.Item(31, 13).Value = 13.330536
$ws.Cells.Item(31, 14).Value = 39.991608
$ws.Cells.Item(31, 15).Value = 0.05365483696904789
$ws.Cells.Item(31, 16).Value = 0.06065723930900103
$ws.Cells.Item(31, 17).Value = 7.04609919596
$ws.Cells.Item(31, 18).Value = 63.41489276364
$ws.Cells.Item(31, 19).Value = 0.0001910351397312933
$ws.Cells.Item(31, 20).Value = 0.0002399899706400255
$ws.Cells.Item(32, 5).Value = 3
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 7).Value = 0.2905836666666666
$ws.Cells.Item(32, 8).Value = 0.8717509999999999
$ws.Cells.Item(32, 9).Value = 0.001957376813162035
$ws.Cells.Item(32, 10).Value = 0.002175106422963506
$ws.Cells.Item(32, 13).Value = 50.6072485
$ws.Cells.Item(32, 14).Value = 101.214497
$ws.Cells.Item(32, 15).Value = 0.2036920096625967
$ws.Cells.Item(32, 16).Value = 0.1535170070198019
$ws.Cells.Item(32, 17).Value = 14.70563982904116
$ws.Cells.Item(32, 18).Value = 88.23383897424699
$ws.Cells.Item(32, 19).Value = 0.000398702016739944
$ws.Cells.Item(32, 20).Value = 0.0003339158280029049
$ws.Cells.Item(33, 5).Value = 3
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = 0.2905836666666666
$ws.Cells.Item(33, 8).Value = 0.8717509999999999
$ws.Cells.Item(33, 9).Value = 0.001957376813162035
$ws.Cells.Item(33, 10).Value = 0.002175106422963506
$ws.Cells.Item(33, 15).Value = 0.5188519664463093
$ws.Cells.Item(33, 16).Value = 0.5865664620849566
$ws.Cells.Item(33, 17).Value = 37.45876019284155
$ws.Cells.Item(33, 18).Value = 337.128841735574
$ws.Cells.Item(33, 19).Value = 0.001015588808585532
$ws.Cells.Item(33, 20).Value = 0.001275844479175969
$ws.Cells.Item(34, 5).Value = 3
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 7).Value = 0.2905836666666666
$ws.Cells.Item(34, 8).Value = 0.8717509999999999
$ws.Cells.Item(34, 9).Value = 0.001957376813162035
$ws.Cells.Item(34, 10).Value = 0.002175106422963506
$ws.Cells.Item(34, 13).Value = 16.762851
$ws.Cells.Item(34, 14).Value = 50.288553
$ws.Cells.Item(34, 15).Value = 0.06746975797083039
$ws.Cells.Item(34, 16).Value = 0.07627512236628199
$ws.Cells.Item(34, 17).Value = 4.871010707367
$ws.Cells.Item(34, 18).Value = 43.839096366303
$ws.Cells.Item(34, 19).Value = 0.0001320637398417578
$ws.Cells.Item(34, 20).Value = 0.0001659065085712273
$ws.Cells.Item(35, 5).Value = 3
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(35, 7).Value = 0.2905836666666666
$ws.Cells.Item(35, 8).Value = 0.8717509999999999
$ws.Cells.Item(35, 9).Value = 0.001957376813162035
$ws.Cells.Item(35, 10).Value = 0.002175106422963506
$ws.Cells.Item(35, 13).Value = 35.4375075
$ws.Cells.Item(35, 14).Value = 70.87501499999999
$ws.Cells.Item(35, 15).Value = 0.1426344512705199
$ws.Cells.Item(35, 16).Value = 0.1074996220678108
$ws.Cells.Item(35, 17).Value = 10.2975608668775
$ws.Cells.Item(35, 18).Value = 61.78536520126499
$ws.Cells.Item(35, 19).Value = 0.0002791893676750059
$ws.Cells.Item(35, 20).Value = 0.0002338231184258447
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 0.2905836666666666
$ws.Cells.Item(36, 8).Value = 0.8717509999999999
$ws.Cells.Item(36, 9).Value = 0.001957376813162035
$ws.Cells.Item(36, 10).Value = 0.002175106422963506
$ws.Cells.Item(36, 13).Value = 3.403012
$ws.Cells.Item(36, 14).Value = 10.209036
$ws.Cells.Item(36, 15).Value = 0.01369697768069593
$ws.Cells.Item(36, 16).Value = 0.01548454715214769
$ws.Cells.Item(36, 17).Value = 0.9888597046706666
$ws.Cells.Item(36, 18).Value = 8.899737342036
$ws.Cells.Item(36, 19).Value = 0.00002681014652259213
$ws.Cells.Item(36, 20).Value = 0.00003368053796731771
$ws.Cells.Item(37, 5).Value = 3
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 0.2905836666666666
$ws.Cells.Item(37, 8).Value = 0.8717509999999999
$ws.Cells.Item(37, 9).Value = 0.001957376813162035
$ws.Cells.Item(37, 10).Value = 0.002175106422963506
$ws.Cells.Item(37, 13).Value = 13.330536
$ws.Cells.Item(37, 14).Value = 39.991608
$ws.Cells.Item(37, 15).Value = 0.05365483696904789
$ws.Cells.Item(37, 16).Value = 0.06065723930900103
$ws.Cells.Item(37, 17).Value = 3.873636029512
$ws.Cells.Item(37, 18).Value = 34.86272426560799
$ws.Cells.Item(37, 19).Value = 0.0001050227337972035
$ws.Cells.Item(37, 20).Value = 0.0001319359508202426
